$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 1 de Junio de 2020 a las 01:05"

# Update country rows: reordered countries (Argentina/Japon/Austria,
# Libia/Guyana/Islas Caimanes/Brunei/Bermudas, Seychelles/Montserrat)
# plus refreshed case statistics for several countries.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value2 = "Estados Unidos"
$ws.Cells.Item(4, 2).Value2 = 1836759
$ws.Cells.Item(4, 3).Value2 = 19939
$ws.Cells.Item(4, 4).Value2 = 541361
$ws.Cells.Item(4, 5).Value2 = 1189222
$ws.Cells.Item(4, 6).Value2 = 0
$ws.Cells.Item(4, 7).Value2 = 619
$ws.Cells.Item(4, 8).Value2 = 106176

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value2 = "Brasil"
$ws.Cells.Item(5, 2).Value2 = 514849
$ws.Cells.Item(5, 3).Value2 = 16409
$ws.Cells.Item(5, 4).Value2 = 206555
$ws.Cells.Item(5, 5).Value2 = 278980
$ws.Cells.Item(5, 6).Value2 = 0
$ws.Cells.Item(5, 7).Value2 = 480
$ws.Cells.Item(5, 8).Value2 = 29314

# Row 34: Colombia
$ws.Cells.Item(34, 1).Value2 = "Colombia"
$ws.Cells.Item(34, 2).Value2 = 29383
$ws.Cells.Item(34, 3).Value2 = 1147
$ws.Cells.Item(34, 4).Value2 = 8543
$ws.Cells.Item(34, 5).Value2 = 19901
$ws.Cells.Item(34, 6).Value2 = 0
$ws.Cells.Item(34, 7).Value2 = 49
$ws.Cells.Item(34, 8).Value2 = 939

# Row 45: Argentina
$ws.Cells.Item(45, 1).Value2 = "Argentina"
$ws.Cells.Item(45, 2).Value2 = 16851
$ws.Cells.Item(45, 3).Value2 = 637
$ws.Cells.Item(45, 4).Value2 = 5336
$ws.Cells.Item(45, 5).Value2 = 10976
$ws.Cells.Item(45, 6).Value2 = 0
$ws.Cells.Item(45, 7).Value2 = 11
$ws.Cells.Item(45, 8).Value2 = 539

# Row 46: Japon
$ws.Cells.Item(46, 1).Value2 = "Japon"
$ws.Cells.Item(46, 2).Value2 = 16851
$ws.Cells.Item(46, 3).Value2 = 47
$ws.Cells.Item(46, 4).Value2 = 14459
$ws.Cells.Item(46, 5).Value2 = 1501
$ws.Cells.Item(46, 6).Value2 = 0
$ws.Cells.Item(46, 7).Value2 = 5
$ws.Cells.Item(46, 8).Value2 = 891

# Row 47: Austria
$ws.Cells.Item(47, 1).Value2 = "Austria"
$ws.Cells.Item(47, 2).Value2 = 16731
$ws.Cells.Item(47, 3).Value2 = 46
$ws.Cells.Item(47, 4).Value2 = 15593
$ws.Cells.Item(47, 5).Value2 = 470
$ws.Cells.Item(47, 6).Value2 = 0
$ws.Cells.Item(47, 7).Value2 = 0
$ws.Cells.Item(47, 8).Value2 = 668

# Row 49: Panama
$ws.Cells.Item(49, 1).Value2 = "Panama"
$ws.Cells.Item(49, 2).Value2 = 13463
$ws.Cells.Item(49, 3).Value2 = 445
$ws.Cells.Item(49, 4).Value2 = 9514
$ws.Cells.Item(49, 5).Value2 = 3613
$ws.Cells.Item(49, 6).Value2 = 0
$ws.Cells.Item(49, 7).Value2 = 6
$ws.Cells.Item(49, 8).Value2 = 336

# Row 60: Chequia
$ws.Cells.Item(60, 1).Value2 = "Chequia"
$ws.Cells.Item(60, 2).Value2 = 9268
$ws.Cells.Item(60, 3).Value2 = 38
$ws.Cells.Item(60, 4).Value2 = 6558
$ws.Cells.Item(60, 5).Value2 = 2390
$ws.Cells.Item(60, 6).Value2 = 0
$ws.Cells.Item(60, 7).Value2 = 1
$ws.Cells.Item(60, 8).Value2 = 320

# Row 165: Libia
$ws.Cells.Item(165, 1).Value2 = "Libia"
$ws.Cells.Item(165, 2).Value2 = 156
$ws.Cells.Item(165, 3).Value2 = 26
$ws.Cells.Item(165, 4).Value2 = 52
$ws.Cells.Item(165, 5).Value2 = 99
$ws.Cells.Item(165, 6).Value2 = 0
$ws.Cells.Item(165, 7).Value2 = 0
$ws.Cells.Item(165, 8).Value2 = 5

# Row 166: Guyana
$ws.Cells.Item(166, 1).Value2 = "Guyana"
$ws.Cells.Item(166, 2).Value2 = 152
$ws.Cells.Item(166, 3).Value2 = 0
$ws.Cells.Item(166, 4).Value2 = 67
$ws.Cells.Item(166, 5).Value2 = 73
$ws.Cells.Item(166, 6).Value2 = 0
$ws.Cells.Item(166, 7).Value2 = 0
$ws.Cells.Item(166, 8).Value2 = 12

# Row 167: Islas Caimanes
$ws.Cells.Item(167, 1).Value2 = "Islas Caimanes"
$ws.Cells.Item(167, 2).Value2 = 141
$ws.Cells.Item(167, 3).Value2 = 0
$ws.Cells.Item(167, 4).Value2 = 68
$ws.Cells.Item(167, 5).Value2 = 72
$ws.Cells.Item(167, 6).Value2 = 0
$ws.Cells.Item(167, 7).Value2 = 0
$ws.Cells.Item(167, 8).Value2 = 1

# Row 168: Brunei
$ws.Cells.Item(168, 1).Value2 = "Brunei"
$ws.Cells.Item(168, 2).Value2 = 141
$ws.Cells.Item(168, 3).Value2 = 0
$ws.Cells.Item(168, 4).Value2 = 138
$ws.Cells.Item(168, 5).Value2 = 1
$ws.Cells.Item(168, 6).Value2 = 0
$ws.Cells.Item(168, 7).Value2 = 0
$ws.Cells.Item(168, 8).Value2 = 2

# Row 169: Bermudas
$ws.Cells.Item(169, 1).Value2 = "Bermudas"
$ws.Cells.Item(169, 2).Value2 = 140
$ws.Cells.Item(169, 3).Value2 = 0
$ws.Cells.Item(169, 4).Value2 = 92
$ws.Cells.Item(169, 5).Value2 = 39
$ws.Cells.Item(169, 6).Value2 = 0
$ws.Cells.Item(169, 7).Value2 = 0
$ws.Cells.Item(169, 8).Value2 = 9

# Row 210: Seychelles
$ws.Cells.Item(210, 1).Value2 = "Seychelles"
$ws.Cells.Item(210, 2).Value2 = 11
$ws.Cells.Item(210, 3).Value2 = 0
$ws.Cells.Item(210, 4).Value2 = 11
$ws.Cells.Item(210, 5).Value2 = 0
$ws.Cells.Item(210, 6).Value2 = 0
$ws.Cells.Item(210, 7).Value2 = 0
$ws.Cells.Item(210, 8).Value2 = 0

# Row 211: Montserrat
$ws.Cells.Item(211, 1).Value2 = "Montserrat"
$ws.Cells.Item(211, 2).Value2 = 11
$ws.Cells.Item(211, 3).Value2 = 0
$ws.Cells.Item(211, 4).Value2 = 10
$ws.Cells.Item(211, 5).Value2 = 0
$ws.Cells.Item(211, 6).Value2 = 0
$ws.Cells.Item(211, 7).Value2 = 0
$ws.Cells.Item(211, 8).Value2 = 0
